{"js": "// Update the report title's building ID (1 -> 4)\nconst titleResults = context.document.body.search(\"Immeuble ID: 1\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\"Immeuble ID: 4\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Update the data row in the table: ID Logement, Type Diagnostic, Date Diagnostic\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.getCell(1, 0).value = \"4\";\n  table.getCell(1, 1).value = \"\u00c9lectricit\u00e9\";\n  table.getCell(1, 2).value = \"2023-08-15\";\n  await context.sync();\n}\n", "ps1": "# Update the report title's building ID (1 -> 4)\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Execute(\"Immeuble ID: 1\", $false, $false, $false, $false, $false, $true, 1, $false, \"Immeuble ID: 4\", 2)\n\n# Update the data row in the table: ID Logement, Type Diagnostic, Date Diagnostic\n$t = $d.Tables.Item(1)\n$t.Cell(2, 1).Range.Text = \"4\"\n$t.Cell(2, 2).Range.Text = \"\u00c9lectricit\u00e9\"\n$t.Cell(2, 3).Range.Text = \"2023-08-15\"\n"}
